$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O3").Value = 1.11
$ws.Range("P3").Value = 6.5
$ws.Range("O6").Value = 1.73
$ws.Range("P6").Value = 2
$ws.Range("S6").Value = 1.75
$ws.Range("T6").Value = 2.05
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 2.4
$ws.Range("J13").Value = 3.1
$ws.Range("K13").Value = 2.27
$ws.Range("L13").Value = 2.87
$ws.Range("N13").Value = 9.25
$ws.Range("S13").Value = 1.29
$ws.Range("T13").Value = 3.3
$ws.Range("V13").Value = 2.65
$ws.Range("Y13").Value = 10
$ws.Range("Z13").Value = 32
$ws.Range("AC13").Value = 9.25
$ws.Range("AD13").Value = 7.4
$ws.Range("AE13").Value = 10.5
$ws.Range("AI13").Value = 16
$ws.Range("AK13").Value = 28
$ws.Range("AL13").Value = 16.5
$ws.Range("AN13").Value = 5
$ws.Range("AO13").Value = 13.5
$ws.Range("AT13").Value = 3.3
$ws.Range("AU13").Value = 6.1
$ws.Range("AW13").Value = 4.8
$ws.Range("AX13").Value = 12
$ws.Range("BC13").Value = 500
$ws.Range("G17").Value = 2.55
$ws.Range("J17").Value = 3.05
$ws.Range("K17").Value = 2.18
$ws.Range("N17").Value = 8
$ws.Range("O17").Value = 1.24
$ws.Range("P17").Value = 3.65
$ws.Range("Q17").Value = 1.72
$ws.Range("R17").Value = 2.05
$ws.Range("S17").Value = 1.36
$ws.Range("T17").Value = 2.9
$ws.Range("U17").Value = 1.6
$ws.Range("V17").Value = 2.22
$ws.Range("W17").Value = 10.25
$ws.Range("X17").Value = 14.5
$ws.Range("Z17").Value = 29
$ws.Range("AA17").Value = 19
$ws.Range("AB17").Value = 24
$ws.Range("AC17").Value = 8
$ws.Range("AE17").Value = 12
$ws.Range("AF17").Value = 45
$ws.Range("AG17").Value = 300
$ws.Range("AM17").Value = 25
$ws.Range("AN17").Value = 4.65
$ws.Range("AO17").Value = 13
$ws.Range("AP17").Value = 19
$ws.Range("AR17").Value = 75
$ws.Range("AS17").Value = 200
$ws.Range("AT17").Value = 2.9
$ws.Range("AU17").Value = 6.7
$ws.Range("AV17").Value = 50
$ws.Range("AW17").Value = 4.6
$ws.Range("G26").Value = 2.05
$ws.Range("H26").Value = 2.7
$ws.Range("G27").Value = 2.57
$ws.Range("I27").Value = 2.7
$ws.Range("G28").Value = 2.85
$ws.Range("I28").Value = 2.45
$ws.Range("G29").Value = 1.96
$ws.Range("H32").Value = 3.8
$ws.Range("I32").Value = 1.8
$ws.Range("M32").Value = 1.02
$ws.Range("N32").Value = 19
$ws.Range("Q32").Value = 1.48
$ws.Range("R32").Value = 2.6
$ws.Range("Y32").Value = 15
$ws.Range("AB32").Value = 29
$ws.Range("AC32").Value = 17
$ws.Range("AD32").Value = 7.5
$ws.Range("AE32").Value = 11
$ws.Range("AH32").Value = 11
$ws.Range("AL32").Value = 13
$ws.Range("AN32").Value = 6
$ws.Range("Q40").Value = 1.8
$ws.Range("R40").Value = 2
$ws.Range("Q41").Value = 2.1
$ws.Range("R41").Value = 1.7
$ws.Range("G60").Value = 4.8
$ws.Range("H60").Value = 3.95
$ws.Range("J60").Value = 4.85
$ws.Range("U60").Value = 1.65
$ws.Range("V60").Value = 2.18
$ws.Range("X60").Value = 24
$ws.Range("AD60").Value = 7
$ws.Range("AN60").Value = 6.7
$ws.Range("M69").Value = 1.04
$ws.Range("N69").Value = 13
$ws.Range("Q69").Value = 1.7
$ws.Range("R69").Value = 2.1
$ws.Range("N75").Value = 8
$ws.Range("Q75").Value = 2.25
$ws.Range("R75").Value = 1.62
$ws.Range("N79").Value = 13
$ws.Range("O79").Value = 1.25
$ws.Range("P79").Value = 3.75
$ws.Range("Q79").Value = 1.83
$ws.Range("R79").Value = 2.03
$ws.Range("G82").Value = 1.95
$ws.Range("I82").Value = 3.9
$ws.Range("J82").Value = 2.63
$ws.Range("L82").Value = 4.5
$ws.Range("O82").Value = 1.36
$ws.Range("P82").Value = 3
$ws.Range("Z82").Value = 17
$ws.Range("AI82").Value = 19
$ws.Range("AJ82").Value = 13
$ws.Range("AO82").Value = 11
$ws.Range("AQ82").Value = 41
$ws.Range("G87").Value = 1.95
$ws.Range("I87").Value = 4.2
$ws.Range("U87").Value = 2.2
$ws.Range("V87").Value = 1.62
$ws.Range("X87").Value = 8
$ws.Range("Z87").Value = 17
$ws.Range("AA87").Value = 21
$ws.Range("AH87").Value = 8.5
$ws.Range("AI87").Value = 19
$ws.Range("AK87").Value = 41
$ws.Range("AO87").Value = 12
$ws.Range("H101").Value = 5.3
$ws.Range("I101").Value = 1.21
$ws.Range("J101").Value = 10.25
$ws.Range("K101").Value = 2.4
$ws.Range("N101").Value = 13.5
$ws.Range("O101").Value = 1.22
$ws.Range("P101").Value = 3.45
$ws.Range("Q101").Value = 1.65
$ws.Range("R101").Value = 1.98
$ws.Range("S101").Value = 1.33
$ws.Range("T101").Value = 3.12
$ws.Range("W101").Value = 29
$ws.Range("X101").Value = 110
$ws.Range("Z101").Value = 600
$ws.Range("AA101").Value = 250
$ws.Range("AC101").Value = 11.25
$ws.Range("AD101").Value = 11.5
$ws.Range("AF101").Value = 200
$ws.Range("AH101").Value = 6
$ws.Range("AK101").Value = 6.3
$ws.Range("AN101").Value = 11.75
$ws.Range("AP101").Value = 80
$ws.Range("AT101").Value = 2.67
$ws.Range("AX101").Value = 5.1
$ws.Range("AZ101").Value = 13
